$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly data row was added to the log. It belongs right before the
# current row 34 (chronologically it slots in among the existing entries),
# so push rows 34:166 down by one to make room, then fill the freed row 34
# with the new observation. Row 167 ends up holding what used to be row 166.
$ws.Rows(34).Insert()

$ws.Range("A34").Value = 8
$ws.Range("B34").Value = "Terminal La Palmera de La Serena"
$ws.Range("C34").Value = "Coquimbo"
$ws.Range("D34").Value = 45145
$ws.Range("E34").Value = 4
$ws.Range("F34").Value = 100114007
$ws.Range("G34").Value = "Jengibre"
$ws.Range("H34").Value = "Sin especificar"
$ws.Range("I34").Value = "Primera"
$ws.Range("J34").Value = 400
$ws.Range("K34").Value = 18000
$ws.Range("L34").Value = 19000
$ws.Range("M34").Value = 18500
$ws.Range("N34").Value = "$/caja 13 kilos"
$ws.Range("O34").Value = "Perú"
$ws.Range("P34").Value = 1423
$ws.Range("Q34").Value = 13
$ws.Range("R34").Value = "Hortaliza"
